$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert column F (ASSISTS) from text-numbers to real numbers for rows 2-41,
# keeping the same displayed values.
$assists = @{
    2 = 0; 3 = 0; 4 = 0; 5 = 0; 6 = 0;
    7 = 1; 8 = 1; 9 = 1; 10 = 1; 11 = 1; 12 = 1; 13 = 1;
    14 = 2; 15 = 2; 16 = 2; 17 = 2;
    18 = 5; 19 = 5; 20 = 5; 21 = 5; 22 = 5; 23 = 5; 24 = 5; 25 = 5;
    26 = 5; 27 = 5; 28 = 5; 29 = 5; 30 = 5; 31 = 5; 32 = 5; 33 = 5; 34 = 5; 35 = 5;
    36 = 6;
    37 = 7; 38 = 7; 39 = 7; 40 = 7; 41 = 7
}

foreach ($row in $assists.Keys) {
    $ws.Cells.Item($row, 6).Value = $assists[$row]
}

# Fix champion names: several rows incorrectly showed other champions
# instead of "Yasuo".
$championRows = @(5, 11, 17, 21, 23, 29, 35, 41)
foreach ($row in $championRows) {
    $ws.Cells.Item($row, 8).Value = "Yasuo"
}
